# UC, add data of minimum duration time
# Adds two new columns ("td1", "td2" - minimum on/off duration for UC) to the
# PV and Slack generator sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# PV sheet: new columns U (td1) and V (td2)
# ---------------------------------------------------------------------------
$pv = $wb.Worksheets.Item("PV")

$pv.Range("U1").Value = "td1"
$pv.Range("V1").Value = "td2"
$pv.Range("U1:V1").Font.Name = $pv.Range("T1").Font.Name
$pv.Range("U1:V1").HorizontalAlignment = -4131   # xlLeft

$pv.Range("U2").Value = 30
$pv.Range("V2").Value = 20
$pv.Range("U3").Value = 45
$pv.Range("V3").Value = 30
$pv.Range("U4").Value = 40
$pv.Range("V4").Value = 23
$pv.Range("U5").Value = 35
$pv.Range("V5").Value = 20
$pv.Range("U2:V5").HorizontalAlignment = -4131   # xlLeft

# ---------------------------------------------------------------------------
# Slack sheet: new columns V (td1) and W (td2)
# ---------------------------------------------------------------------------
$slack = $wb.Worksheets.Item("Slack")

$slack.Range("V1").Value = "td1"
$slack.Range("W1").Value = "td2"
$slack.Range("V1:W1").Font.Name = $slack.Range("U1").Font.Name
$slack.Range("V1:W1").HorizontalAlignment = -4131  # xlLeft
$slack.Range("V1:W1").VerticalAlignment = -4160    # xlTop

$slack.Range("V2").Value = 50
$slack.Range("W2").Value = 30
$slack.Range("V2:W2").HorizontalAlignment = -4131  # xlLeft
$slack.Range("V2:W2").VerticalAlignment = -4160    # xlTop
